# Auto-generated: apply scheduled price-data refresh to Lamia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 721.8889
$ws.Range("I58").Value = 721.8889
$ws.Range("K58").Value = 2165.6667
$ws.Range("M58").Value = -2015.6667

# Row 112
$ws.Range("H112").Value = 1526.5
$ws.Range("J112").Value = 1526.5
$ws.Range("L112").Value = 4579.5
$ws.Range("N112").Value = -6795.5

# Row 138
$ws.Range("H138").Value = 5341.36
$ws.Range("J138").Value = 6068.0586
$ws.Range("L138").Value = 18204.1758
$ws.Range("N138").Value = -28484.1758

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 2090.3635
$ws.Range("I22").Value = 1583.3334
$ws.Range("K22").Value = 1583.3334
$ws.Range("M22").Value = -1284.3334

# Row 32
$ws.Range("H32").Value = 1776.1731
$ws.Range("I32").Value = 1657.8636
$ws.Range("J32").Value = 2426.875
$ws.Range("K32").Value = 1657.8636
$ws.Range("L32").Value = 2426.875
$ws.Range("M32").Value = -1370.8636
$ws.Range("N32").Value = -3000.875

# Row 61
$ws.Range("H61").Value = 5264.8276
$ws.Range("I61").Value = 4131.5
$ws.Range("K61").Value = 4131.5
$ws.Range("M61").Value = -3919.5

# Row 97
$ws.Range("H97").Value = 952.8077
$ws.Range("I97").Value = 819.2381
$ws.Range("K97").Value = 819.2381
$ws.Range("M97").Value = -323.2381

# Row 136
$ws.Range("H136").Value = 5264.8276
$ws.Range("I136").Value = 4131.5
$ws.Range("K136").Value = 12394.5
$ws.Range("M136").Value = -9844.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5617.6665
$ws.Range("J20").Value = 9091.200000000001
$ws.Range("L20").Value = 9091.200000000001
$ws.Range("N20").Value = -9585.200000000001

# Row 105
$ws.Range("H105").Value = 17141.15
$ws.Range("I105").Value = 16362.5625
$ws.Range("J105").Value = 20255.5
$ws.Range("K105").Value = 16362.5625
$ws.Range("L105").Value = 20255.5
$ws.Range("M105").Value = -14615.5625
$ws.Range("N105").Value = -23749.5

# Row 134
$ws.Range("H134").Value = 3731.1892
$ws.Range("I134").Value = 2195.6667
$ws.Range("K134").Value = 6587.000100000001
$ws.Range("M134").Value = -4052.000100000001

# Row 141
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1864.0588
$ws.Range("I16").Value = 1370.7142
$ws.Range("K16").Value = 1370.7142
$ws.Range("M16").Value = -1083.7142

# Row 22
$ws.Range("H22").Value = 2854.5715
$ws.Range("I22").Value = 1555.2222
$ws.Range("K22").Value = 1555.2222
$ws.Range("M22").Value = -1205.2222

# Row 58
$ws.Range("H58").Value = 4449.5854
$ws.Range("I58").Value = 3182.577
$ws.Range("K58").Value = 3182.577
$ws.Range("M58").Value = -2979.577

# Row 62
$ws.Range("H62").Value = 16625.5
$ws.Range("I62").Value = 19333.334
$ws.Range("K62").Value = 19333.334
$ws.Range("M62").Value = -18709.334

# Row 65
$ws.Range("H65").Value = 16625.5
$ws.Range("I65").Value = 19333.334
$ws.Range("K65").Value = 96666.67
$ws.Range("M65").Value = -93546.67

# Row 113
$ws.Range("H113").Value = 1864.0588
$ws.Range("I113").Value = 1370.7142
$ws.Range("K113").Value = 1370.7142
$ws.Range("M113").Value = 799.2858000000001

# Row 122
$ws.Range("H122").Value = 5315.089
$ws.Range("I122").Value = 3298.3872
$ws.Range("J122").Value = 9780.643
$ws.Range("K122").Value = 9895.161599999999
$ws.Range("L122").Value = 29341.929
$ws.Range("M122").Value = -7445.161599999999
$ws.Range("N122").Value = -34241.929

# Row 136
$ws.Range("H136").Value = 4449.5854
$ws.Range("I136").Value = 3182.577
$ws.Range("K136").Value = 9547.731
$ws.Range("M136").Value = -6997.731

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 447.8
$ws.Range("J23").Value = 486.44446
$ws.Range("L23").Value = 1459.33338
$ws.Range("N23").Value = -1929.33338

# Row 107
$ws.Range("H107").Value = 1450.9667
$ws.Range("I107").Value = 1463.9231
$ws.Range("J107").Value = 1441.0588
$ws.Range("K107").Value = 4391.7693
$ws.Range("L107").Value = 4323.1764
$ws.Range("M107").Value = -2471.7693
$ws.Range("N107").Value = -8163.1764

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1396.95
$ws.Range("I97").Value = 1537.0834
$ws.Range("J97").Value = 1186.75
$ws.Range("K97").Value = 1537.0834
$ws.Range("L97").Value = 1186.75
$ws.Range("M97").Value = -1041.0834
$ws.Range("N97").Value = -2178.75

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 1003.3333
$ws.Range("J5").Value = 1003.3333
$ws.Range("L5").Value = 1003.3333
$ws.Range("N5").Value = -1229.3333

# Row 7
$ws.Range("H7").Value = 3947.1667
$ws.Range("I7").Value = 3947.1667
$ws.Range("K7").Value = 3947.1667
$ws.Range("M7").Value = -3835.1667

# Row 22
$ws.Range("H22").Value = 2510.375
$ws.Range("I22").Value = 1733.3334
$ws.Range("J22").Value = 2976.6
$ws.Range("K22").Value = 1733.3334
$ws.Range("L22").Value = 2976.6
$ws.Range("M22").Value = -1438.3334
$ws.Range("N22").Value = -3566.6

# Row 27
$ws.Range("H27").Value = 2510.375
$ws.Range("I27").Value = 1733.3334
$ws.Range("J27").Value = 2976.6
$ws.Range("K27").Value = 1733.3334
$ws.Range("L27").Value = 2976.6
$ws.Range("M27").Value = -1626.3334
$ws.Range("N27").Value = -3190.6

# Row 46
$ws.Range("H46").Value = 3403.9092
$ws.Range("I46").Value = 1221
$ws.Range("J46").Value = 3889
$ws.Range("K46").Value = 1221
$ws.Range("L46").Value = 3889
$ws.Range("M46").Value = -1033
$ws.Range("N46").Value = -4265

# Row 55
$ws.Range("H55").Value = 8334969.5
$ws.Range("I55").Value = 16666926
$ws.Range("K55").Value = 16666926
$ws.Range("M55").Value = -16666753

# Row 82
$ws.Range("H82").Value = 7573.263
$ws.Range("I82").Value = 5439.3
$ws.Range("K82").Value = 5439.3
$ws.Range("M82").Value = -5078.3

# Row 85
$ws.Range("H85").Value = 7573.263
$ws.Range("I85").Value = 5439.3
$ws.Range("K85").Value = 5439.3
$ws.Range("M85").Value = -4191.3

# Row 93
$ws.Range("H93").Value = 11272.489
$ws.Range("I93").Value = 6007.9
$ws.Range("K93").Value = 6007.9
$ws.Range("M93").Value = -4759.9

# Row 126
$ws.Range("H126").Value = 3947.1667
$ws.Range("I126").Value = 3947.1667
$ws.Range("K126").Value = 11841.5001
$ws.Range("M126").Value = -9371.500100000001

# Row 136
$ws.Range("H136").Value = 7036.927
$ws.Range("I136").Value = 4775.161
$ws.Range("J136").Value = 14048.4
$ws.Range("K136").Value = 14325.483
$ws.Range("L136").Value = 42145.2
$ws.Range("M136").Value = -11775.483
$ws.Range("N136").Value = -47245.2

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4378.3335
$ws.Range("I126").Value = 3592.75
$ws.Range("J126").Value = 5006.8
$ws.Range("K126").Value = 10778.25
$ws.Range("L126").Value = 15020.4
$ws.Range("M126").Value = -8308.25
$ws.Range("N126").Value = -19960.4

# Row 136
$ws.Range("H136").Value = 3380.2285
$ws.Range("I136").Value = 1800.3334
$ws.Range("K136").Value = 5401.0002
$ws.Range("M136").Value = -2851.0002
